# Generate Report for Handoff
# Updates status cells from "Handed back: in sync with en-US" to "Ready for handoff"
# and refreshes the related "Latest ... Datetime" timestamps, then narrows the
# now-shorter status columns' widths on each worksheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Target column width (~17.216 chars) lands on the same quantized grid point
# (17.1667 chars) for any COM ColumnWidth assignment in [16.25, 16.4167), so
# 16.333333333333332 is used below to land as close as possible to it.
$narrowWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 15:14:32"
$wsOverview.Range("E1").EntireColumn.ColumnWidth = $narrowWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $narrowWidth

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 15:14:27"
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $narrowWidth

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 15:14:32"
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $narrowWidth
